$d = $word.ActiveDocument

# Replace the leading portion of the title text while keeping the
# trailing " Web-App" in place.
$d.Content.Find.Execute("D7 Auto Service Center", $true, $false, $false, $false, $false, $true, 1, $false, "Barangay South Signal Village", 2)

# The remaining " Web-App" text needs to end up in its own run (same
# formatting as before) instead of being merged into the preceding run.
# Locate it and force a run boundary by nudging a character-formatting
# property away from, and then back to, its original value.
$rng = $d.Content
$rng.Find.Execute(" Web-App", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Size = 99
$rng.Font.Size = 22
